$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: Coin (B), Link (C), Price (D), Volume 1h (E)
$rowUpdates = @(
    @{ Row = 2; D = '23.323.96'; E = '  -0.42%  ' },
    @{ Row = 3; D = '1.625.15'; E = '  -0.81%  ' },
    @{ Row = 4; E = '  +0.03%  ' },
    @{ Row = 5; E = '  -0.06%  ' },
    @{ Row = 6; D = '302.83'; E = '  -0.65%  ' },
    @{ Row = 7; E = '  +0.53%  ' },
    @{ Row = 8; D = '0.3619'; E = '  +0.21%  ' },
    @{ Row = 9; D = '51.24'; E = '  -1.27%  ' },
    @{ Row = 10; D = '0.08141'; E = '  +0.41%  ' },
    @{ Row = 11; D = '1.224'; E = '  -2.00%  ' },
    @{ Row = 12; E = '  -0.14%  ' },
    @{ Row = 13; D = '22.29'; E = '  -1.89%  ' },
    @{ Row = 14; D = '6.476'; E = '  -1.42%  ' },
    @{ Row = 15; D = '0.00001239'; E = '  -1.98%  ' },
    @{ Row = 16; D = '7.270'; E = '  +0.16%  ' },
    @{ Row = 17; D = '1.622.26'; E = '  -0.77%  ' },
    @{ Row = 18; D = '93.92'; E = '  -0.19%  ' },
    @{ Row = 19; D = '0.06935'; E = '  +0.69%  ' },
    @{ Row = 20; D = '17.51'; E = '  -3.02%  ' },
    @{ Row = 21; D = '6.524' },
    @{ Row = 22; E = '  +0.01%  ' },
    @{ Row = 23; D = '12.51'; E = '  -1.41%  ' },
    @{ Row = 24; D = '23.330.10'; E = '  -0.44%  ' },
    @{ Row = 25; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '3.103'; E = '  +2.37%  ' },
    @{ Row = 26; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '2.463'; E = '  +2.22%  ' },
    @{ Row = 27; D = '21.17'; E = '  +0.02%  ' },
    @{ Row = 28; D = '150.47'; E = '  -0.74%  ' },
    @{ Row = 29; D = '5.266'; E = '  -0.81%  ' },
    @{ Row = 30; D = '132.74'; E = '  -1.92%  ' },
    @{ Row = 31; D = '1.796.51'; E = '  -1.02%  ' },
    @{ Row = 32; D = '6.719'; E = '  -0.15%  ' },
    @{ Row = 33; D = '2.172'; E = '  -4.65%  ' },
    @{ Row = 34; D = '1.061'; E = '  +11.71%  ' },
    @{ Row = 35; D = '11.25'; E = '  +9.74%  ' },
    @{ Row = 36; D = '0.02757'; E = '  -1.54%  ' },
    @{ Row = 37; D = '0.08781'; E = '  +0.27%  ' },
    @{ Row = 38; D = '0.2478'; E = '  -1.27%  ' },
    @{ Row = 39; D = '0.07104'; E = '  -1.39%  ' },
    @{ Row = 40; D = '5.996'; E = '  -0.61%  ' },
    @{ Row = 41; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.6976'; E = '  -0.63%  ' },
    @{ Row = 42; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.338'; E = '  -2.34%  ' },
    @{ Row = 43; D = '16.00'; E = '  -0.16%  ' },
    @{ Row = 44; D = '12.06'; E = '  -2.69%  ' },
    @{ Row = 45; D = '0.6466'; E = '  -0.08%  ' },
    @{ Row = 46; D = '1.000'; E = '  -0.05%  ' },
    @{ Row = 47; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '2.269'; E = '  -1.91%  ' },
    @{ Row = 48; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '3.959'; E = '  -1.01%  ' },
    @{ Row = 49; D = '0.07967'; E = '  +0.03%  ' },
    @{ Row = 50; E = '  -1.65%  ' },
    @{ Row = 51; E = '  -1.15%  ' }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        # Column D holds price strings that look numeric (e.g. "302.83").
        # Force text so Excel doesn't auto-convert them to numbers.
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
